# Add two new localization rows (app.urlcleaner.confirmBtn / app.urlcleaner.ruleTitle)
# to the end of the i18n table on Sheet1, extending the sheet from row 53 to row 55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles, row height, column layout) from the last existing
# data row (53) down into the two new rows (54:55) so the new rows match the
# look of the rest of the table.
$ws.Range("A53:E53").Copy()
$ws.Range("A54:E55").PasteSpecial(-4122)
$ws.Rows.Item(54).RowHeight = 20.1
$ws.Rows.Item(55).RowHeight = 20.1

# Row 54 - app.urlcleaner.confirmBtn
$ws.Range("A54").Value = "app.urlcleaner.confirmBtn"
$ws.Range("B54").Value = "净化"
$ws.Range("C54").Value = "Clean URL"

# Row 55 - app.urlcleaner.ruleTitle
$ws.Range("A55").Value = "app.urlcleaner.ruleTitle"
$ws.Range("B55").Value = "规则"
$ws.Range("C55").Value = "Rules"
